# 煤油.xlsx edit script
#
# Logical change (per commit diff):
#   1. Within every 4-row year block (rows 2-5, 6-9, 10-13, ... 62-64), the
#      second row ("B" sub-period) and third row ("C" sub-period) swap their
#      entire data content (columns A-E), while staying in the same physical
#      row position. The last block only has 3 rows (62-64) but still swaps
#      its B/C rows (63 & 64).
#   2. Columns F and G (煤油产销率 / 煤油销售量 - the single-period, as
#      opposed to cumulative, figures) are removed entirely, shrinking the
#      used range from A1:G64 down to A1:E64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (B-row, C-row) that swap within each 4-row year block.
$pairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16), @(19,20), @(23,24), @(27,28), @(31,32),
    @(35,36), @(39,40), @(43,44), @(47,48), @(51,52), @(55,56), @(59,60), @(63,64)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]

    # Capture the current content of both rows (columns A-E) up front so the
    # writes below don't clobber a value we still need to read.
    $a1 = $ws.Range("A$r1").Value2
    $b1 = $ws.Range("B$r1").Value2
    $c1 = $ws.Range("C$r1").Value2
    $d1 = $ws.Range("D$r1").Value2
    $e1 = $ws.Range("E$r1").Value2

    $a2 = $ws.Range("A$r2").Value2
    $b2 = $ws.Range("B$r2").Value2
    $c2 = $ws.Range("C$r2").Value2
    $d2 = $ws.Range("D$r2").Value2
    $e2 = $ws.Range("E$r2").Value2

    $ws.Range("A$r1").Value2 = $a2
    $ws.Range("B$r1").Value2 = $b2
    $ws.Range("D$r1").Value2 = $d2
    $ws.Range("E$r1").Value2 = $e2

    $ws.Range("A$r2").Value2 = $a1
    $ws.Range("B$r2").Value2 = $b1
    $ws.Range("D$r2").Value2 = $d1
    $ws.Range("E$r2").Value2 = $e1

    # Column C needs special handling: if both cells are blank, leave them
    # alone so the underlying (originally inline, now shared-string) empty
    # text cell is preserved rather than writing "" and deleting the cell
    # outright.
    if (-not ($c1 -eq "" -and $c2 -eq "")) {
        $ws.Range("C$r1").Value2 = $c2
        $ws.Range("C$r2").Value2 = $c1
    }
}

# Drop the 煤油产销率 / 煤油销售量 columns (F, G) - header + all data rows.
$ws.Range("F1:G64").Delete()
